# Add a new "31-jul" data column (AO) to the right of the existing "28-jul"
# column (AN) on Sheet1, carrying one value per product row (2-18).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header
$ws.Range("AO1").Value = "31-jul"

# Data rows (row -> value)
$values = @{
    2  = 0
    3  = 17.434819783586502
    4  = 16.500922457349205
    5  = 20.153971215047854
    6  = 0
    7  = 8.4121172344980533
    8  = 6.5614872419778925
    9  = 15.648249761236167
    10 = 15.534729889681385
    11 = 7.9180654772955563
    12 = 0
    13 = 11.864538843514532
    14 = 0
    15 = 0
    16 = 0.22576509844920953
    17 = 0
    18 = 0
}

foreach ($row in $values.Keys) {
    $ws.Cells.Item($row, 41).Value = $values[$row]
}

# Match the author's final selection/active cell captured in the saved file.
$ws.Range("AQ7").Select()
